$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Rows(7).Insert()
$ws.Range("A7").Value2 = "patient_id"
$ws.Range("B7").Value2 = "VARCHR(16)"
$ws.Range("C7").Value2 = "N"
$ws.Range("D7").Value2 = "환자코드 - 시연에서 P0001"
